$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Write the text (A-column / header) values first, in the exact order
# needed so the workbook's shared-string table is rebuilt with the same
# ordering as the target file: surviving strings first (in their original
# relative order), followed by the brand-new product names in the order
# they were introduced. ---

# Surviving strings (content unchanged, just re-indexed by the engine)
$ws.Range("A9").Value = "Уголь древесно-брикетный (евроуголь) 10 кг"
$ws.Range("A10").Value = "Уголь брикетный 10 кг в коробке"
$ws.Range("B1").Value = "Наличный расчет"
$ws.Range("C1").Value = "Без НДС"
$ws.Range("D1").Value = "С НДС"
$ws.Range("A1").Value = "Наименование"

# Brand-new product name strings
$ws.Range("A7").Value = "Уголь древесный 10 л"
$ws.Range("A8").Value = "Уголь древесный 20 л"
$ws.Range("A2").Value = "Уголь березовый 10 литров"
$ws.Range("A3").Value = "Уголь березовый 20 литров"
$ws.Range("A5").Value = "Уголь березовый 10 литров (эконом)"
$ws.Range("A6").Value = "Древесно-угольный микс 3кг"
$ws.Range("A4").Value = "Уголь березовый 50л (10кг)"

# --- Update the price columns (B, C, D) for every row ---

# Row 2 - Уголь березовый 10 литров
$ws.Range("B2").Value = 130
$ws.Range("C2").Value = 143
$ws.Range("D2").Value = 175.2

# Row 3 - Уголь березовый 20 литров
$ws.Range("B3").Value = 260
$ws.Range("C3").Value = 286
$ws.Range("D3").Value = 343.2

# Row 4 - Уголь березовый 50л (10кг)
$ws.Range("B4").Value = 500
$ws.Range("C4").Value = 550
$ws.Range("D4").Value = 660

# Row 5 - Уголь березовый 10 литров (эконом)
$ws.Range("B5").Value = 120
$ws.Range("C5").Value = 132
$ws.Range("D5").Value = 160.8

# Row 6 - Древесно-угольный микс 3кг
$ws.Range("B6").Value = 200
$ws.Range("C6").Value = 220
$ws.Range("D6").Value = 264

# Row 7 - Уголь древесный 10 л
$ws.Range("B7").Value = 55
$ws.Range("C7").Value = 60.5
$ws.Range("D7").Value = 72.6

# Row 8 - Уголь древесный 20 л
$ws.Range("B8").Value = 110
$ws.Range("C8").Value = 121
$ws.Range("D8").Value = 145.2

# Row 9 - Уголь древесно-брикетный (евроуголь) 10 кг
$ws.Range("B9").Value = 500
$ws.Range("C9").Value = 550
$ws.Range("D9").Value = 660

# Row 10 - Уголь брикетный 10 кг в коробке
$ws.Range("B10").Value = 450
$ws.Range("C10").Value = 495
$ws.Range("D10").Value = 594

# Update the active selection to match the edited workbook
$ws.Range("D10").Select()
